$d = $word.ActiveDocument

# 1) Rename the "Recommend Archaeology Condition to be Satisfied" heading.
$d.Content.Find.Execute("Recommend Archaeology Condition to be Satisfied", $true, $false, $false, $false, $false, $true, 1, $false, "Recommend Archaeology Condition(s)", 2)

# 2) Move the hidden "_GoBack" bookmark from just before "<Proposal Description>"
#    down into the GLAAS paragraph, right between "...follows the " and
#    "National Planning Policy Framework...".
$r = $d.Content
$r.Find.Execute("Our advice follows the ")
$splitPoint = $r.End

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
